$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 57 ---
$ws.Range("B57").Value = "Forgot sign in details? Not to worry. Just tap the relevant button below to fix this."
$ws.Range("A57").Value = "FORGOT_SIGN_IN_DETAILS_LABEL_TXT"
$ws.Range("E57").Value = "device-accessibilityid"

# --- Row 58 ---
$ws.Range("B58").Value = "Forgot username"
$ws.Range("A58").Value = "FORGOT_USER_NAME_BTN"
$ws.Range("E58").Value = "device-accessibilityid"

# --- Row 59 ---
$ws.Range("B59").Value = "Forgot password"
$ws.Range("A59").Value = "FORGOT_PWD_BTN"
$ws.Range("E59").Value = "device-accessibilityid"

# --- Row 60 ---
$ws.Range("B60").Value = "Cancel"
$ws.Range("A60").Value = "CANCEL_BTN"
$ws.Range("E60").Value = "device-accessibilityid"

# --- Row 61 ---
$ws.Range("B61").Value = "Accept All"
$ws.Range("A61").Value = "ACCEPT_ALL_WEBSITE_BTN"
$ws.Range("E61").Value = "device-accessibilityid"

# --- Rows 62-64 (DD / MM / YYYY website fields) ---
$ws.Range("B64").Value = "//XCUIElementTypeTextField[@value='YYYY']"
$ws.Range("B62").Value = "//XCUIElementTypeTextField[@value='DD']"
$ws.Range("A62").Value = "DD_FIELD_IN_WEBSITE"
$ws.Range("B63").Value = "//XCUIElementTypeTextField[@value='MM']"
$ws.Range("A63").Value = "MM_FIELD_IN_WEBSITE"
$ws.Range("A64").Value = "YYYY_FIELD_IN_WEBSITE"
$ws.Range("E62").Value = "device-xpath"
$ws.Range("E63").Value = "device-xpath"
$ws.Range("E64").Value = "device-xpath"

# --- Rows 65-67 (account number / DOB / account number static text) ---
$ws.Range("B65").Value = '(//XCUIElementTypeOther[@name="Forgotten Username"])[1]/XCUIElementTypeTextField[4]'
$ws.Range("A65").Value = "ACCOUNT_NUMBER_FIELD_IN_WEBSITE"
$ws.Range("B66").Value = "//XCUIElementTypeStaticText[@name='Date Of Birth']"
$ws.Range("A66").Value = "DOB_STATIC_TEXT_IN_WEBSITE"
$ws.Range("A67").Value = "ACCOUNT_NUMBER_STATIC_TEXT_IN_WEBSITE"
$ws.Range("B67").Value = "//XCUIElementTypeStaticText[@name='Account Number']"
$ws.Range("E65").Value = "device-xpath"
$ws.Range("E66").Value = "device-xpath"
$ws.Range("E67").Value = "device-xpath"

# --- Row 68 ---
$ws.Range("B68").Value = "//XCUIElementTypeStaticText[@name='Forgotten Password']"
$ws.Range("A68").Value = "FORGOTTEN_PWD_TXT_IN_WEBSITE"
$ws.Range("E68").Value = "device-xpath"

# --- Row 69 ---
$ws.Range("B69").Value = "//XCUIElementTypeStaticText[@name='Username']"
$ws.Range("A69").Value = "USERNAME_TXT_IN_WEBSITE"
$ws.Range("E69").Value = "device-xpath"

# --- Row 70 ---
$ws.Range("B70").Value = "//XCUIElementTypeButton[@name='Next']"
$ws.Range("A70").Value = "NEXT_BTN_IN_WEBSITE"
$ws.Range("E70").Value = "device-xpath"

# --- Row 71 ---
$ws.Range("B71").Value = "//XCUIElementTypeStaticText[@name='Forgotten Username']"
$ws.Range("A71").Value = "FORGOTTEN_USERNAME_TXT_IN_WEBSITE"
$ws.Range("E71").Value = "device-xpath"

# Update the view state to match where the editor left off.
$excel.ActiveWindow.ScrollRow = 39
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B72").Select()
